# Rows 11-14 and 17-19 of the "Artfynd" sheet got their per-observation
# data (Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn,
# Auktor, Ost, Nord and the public-comment text) shuffled between rows,
# while the shared/location columns (locality, date, county, ...) stayed
# put. Concretely:
#   row 11 <-> row 12   (full swap)
#   row 13 <-> row 14   (full swap)
#   row 17 -> row 19 -> row 18 -> row 17   (3-way cyclic rotation:
#       new row17 = old row19, new row18 = old row17, new row19 = old row18)
#
# Apply the change by writing the post-edit values directly into the
# affected cells (A, B, E, F, G, H, Q, R, AC) for each row.
# NOTE: named parameters (-Row 11) do not bind in this PS host, so the
# helper below takes plain positional arguments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ObsRow {
    param($Row, $A, $B, $E, $F, $G, $H, $Q, $R, $AC)

    $ws.Range("A$Row").Value = $A
    $ws.Range("B$Row").Value = $B
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G
    $ws.Range("H$Row").Value = $H
    $ws.Range("Q$Row").Value = $Q
    $ws.Range("R$Row").Value = $R
    if ($AC -eq $null) {
        $ws.Range("AC$Row").Value = ""
    } else {
        $ws.Range("AC$Row").Value = $AC
    }
}

# Row 11 <= old row 12
Set-ObsRow 11 131064766 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 442271 7039174 "Ringhack äldre"

# Row 12 <= old row 11
Set-ObsRow 12 131064783 91828 5432 "Granticka" "Porodaedalea chrysoloma s.lat." "" 442292 7039182 $null

# Row 13 <= old row 14
Set-ObsRow 13 131064779 91804 1108 "Harticka" "Pelloporus leporinus" "(Fr.) Krieglst." 442245 7039149 $null

# Row 14 <= old row 13
Set-ObsRow 14 131064763 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 442230 7039147 "Ringhack äldre"

# Row 17 <= old row 19
Set-ObsRow 17 131064780 91804 1108 "Harticka" "Pelloporus leporinus" "(Fr.) Krieglst." 442259 7039181 $null

# Row 18 <= old row 17
Set-ObsRow 18 131064772 57884 100109 "Tretåig hackspett" "Picoides tridactylus" "(Linnaeus, 1758)" 442099 7039220 "Bohål ca 3m upp i grantickerötad granhögstubbe Även ett påbörjat på 2m"

# Row 19 <= old row 18
Set-ObsRow 19 131064781 91804 1108 "Harticka" "Pelloporus leporinus" "(Fr.) Krieglst." 442200 7039150 $null
